# Generate Report for Handoff
# Inserts a new "234eb279-eab2-4444-ab51-ed33d460a200" file row above the
# existing "491341b3-b61d-4470-859f-90f13b1eafa5" row on each of the three
# sheets (Overview, zh-cn, de-de), marking it "Ready for handoff" / "Include".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Hyperlinks don't move with inserted rows, so drop them all and rebuild
# after the data is in its final place.
$ws.Range("A1").Hyperlinks.Delete()

$ws.Rows(2).Insert()

$ws.Range("A2").Value = "234eb279-eab2-4444-ab51-ed33d460a200.md"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/447d71adbda634f1fe8589a49941a213a98042b3/e2e/234eb279-eab2-4444-ab51-ed33d460a200.md", "", "", "234eb279-eab2-4444-ab51-ed33d460a200.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/447d71adbda634f1fe8589a49941a213a98042b3/e2e/491341b3-b61d-4470-859f-90f13b1eafa5.md", "", "", "491341b3-b61d-4470-859f-90f13b1eafa5.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/447d71adbda634f1fe8589a49941a213a98042b3/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A1").Hyperlinks.Delete()

$ws.Rows(2).Insert()

$ws.Range("A2").Value = "234eb279-eab2-4444-ab51-ed33d460a200.md"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "234eb279-eab2-4444-ab51-ed33d460a200.6ed28c1fb8679de01dcd00b61946f3795714dead.zh-cn.xlf"
$ws.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("D2").Value = "2016-03-08 14:28:15"
$ws.Range("G2").Value = "0001-01-01 00:00:00"
$ws.Range("H2").Value = "Include"

$ws.Range("D3").Value = "2016-03-08 14:27:49"
$ws.Range("H3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/447d71adbda634f1fe8589a49941a213a98042b3/e2e/234eb279-eab2-4444-ab51-ed33d460a200.md", "", "", "234eb279-eab2-4444-ab51-ed33d460a200.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/07e96a8684fdb24a38861919f9e14a0f2ba28503/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/234eb279-eab2-4444-ab51-ed33d460a200.6ed28c1fb8679de01dcd00b61946f3795714dead.zh-cn.xlf", "", "", "234eb279-eab2-4444-ab51-ed33d460a200.6ed28c1fb8679de01dcd00b61946f3795714dead.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/447d71adbda634f1fe8589a49941a213a98042b3/e2e/491341b3-b61d-4470-859f-90f13b1eafa5.md", "", "", "491341b3-b61d-4470-859f-90f13b1eafa5.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/07e96a8684fdb24a38861919f9e14a0f2ba28503/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/491341b3-b61d-4470-859f-90f13b1eafa5.64769047756e11626bcd77981eff3916691eed7b.zh-cn.xlf", "", "", "491341b3-b61d-4470-859f-90f13b1eafa5.64769047756e11626bcd77981eff3916691eed7b.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/447d71adbda634f1fe8589a49941a213a98042b3/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A1").Hyperlinks.Delete()

$ws.Rows(2).Insert()

$ws.Range("A2").Value = "234eb279-eab2-4444-ab51-ed33d460a200.md"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "234eb279-eab2-4444-ab51-ed33d460a200.6ed28c1fb8679de01dcd00b61946f3795714dead.de-de.xlf"
$ws.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("D2").Value = "2016-03-08 14:28:20"
$ws.Range("G2").Value = "0001-01-01 00:00:00"
$ws.Range("H2").Value = "Include"

$ws.Range("D3").Value = "2016-03-08 14:27:54"
$ws.Range("H3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/447d71adbda634f1fe8589a49941a213a98042b3/e2e/234eb279-eab2-4444-ab51-ed33d460a200.md", "", "", "234eb279-eab2-4444-ab51-ed33d460a200.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a4107492ba6577a68b0b37225eba5347afcd98cd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/234eb279-eab2-4444-ab51-ed33d460a200.6ed28c1fb8679de01dcd00b61946f3795714dead.de-de.xlf", "", "", "234eb279-eab2-4444-ab51-ed33d460a200.6ed28c1fb8679de01dcd00b61946f3795714dead.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/447d71adbda634f1fe8589a49941a213a98042b3/e2e/491341b3-b61d-4470-859f-90f13b1eafa5.md", "", "", "491341b3-b61d-4470-859f-90f13b1eafa5.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a4107492ba6577a68b0b37225eba5347afcd98cd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/491341b3-b61d-4470-859f-90f13b1eafa5.64769047756e11626bcd77981eff3916691eed7b.de-de.xlf", "", "", "491341b3-b61d-4470-859f-90f13b1eafa5.64769047756e11626bcd77981eff3916691eed7b.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/447d71adbda634f1fe8589a49941a213a98042b3/.localization-config", "", "", ".localization-config")
